$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Halve every numeric value in the data range B2:K6 (mean/sd columns for
# the Big Five cluster statistics), leaving the Cluster id column (A) untouched.
$range = $ws.Range("B2:K6")
for ($r = 1; $r -le $range.Rows.Count; $r++) {
    for ($c = 1; $c -le $range.Columns.Count; $c++) {
        $cell = $range.Cells.Item($r, $c)
        $cell.Value2 = $cell.Value2 / 2
    }
}
